$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.862146333333333
$ws.Range("H2").Value = 8.586439
$ws.Range("I2").Value = 0.9778268096017091
$ws.Range("J2").Value = 0.9778268096017091
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.677453
$ws.Range("N2").Value = 5.032359
$ws.Range("O2").Value = 0.2543757648546719
$ws.Range("P2").Value = 0.2543757648546719
$ws.Range("Q2").Value = 4.801115953289
$ws.Range("R2").Value = 43.210043579601
$ws.Range("S2").Value = 0.2487354425878384
$ws.Range("T2").Value = 0.2487354425878384

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.862146333333333
$ws.Range("H3").Value = 8.586439
$ws.Range("I3").Value = 0.9778268096017091
$ws.Range("J3").Value = 0.9778268096017091
$ws.Range("O3").Value = 0.5611569328879042
$ws.Range("P3").Value = 0.5611569328879042
$ws.Range("Q3").Value = 10.59133720669522
$ws.Range("R3").Value = 95.32203486025699
$ws.Range("S3").Value = 0.5487142933716598
$ws.Range("T3").Value = 0.5487142933716598

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.862146333333333
$ws.Range("H4").Value = 8.586439
$ws.Range("I4").Value = 0.9778268096017091
$ws.Range("J4").Value = 0.9778268096017091
$ws.Range("M4").Value = 1.148238333333333
$ws.Range("N4").Value = 3.444715
$ws.Range("O4").Value = 0.174123510034034
$ws.Range("P4").Value = 0.174123510034034
$ws.Range("Q4").Value = 3.286426135542778
$ws.Range("R4").Value = 29.577835219885
$ws.Range("S4").Value = 0.1702626362932307
$ws.Range("T4").Value = 0.1702626362932306

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.862146333333333
$ws.Range("H5").Value = 8.586439
$ws.Range("I5").Value = 0.9778268096017091
$ws.Range("J5").Value = 0.9778268096017091
$ws.Range("M5").Value = 0.06821100000000001
$ws.Range("N5").Value = 0.204633
$ws.Range("O5").Value = 0.01034379222338988
$ws.Range("P5").Value = 0.01034379222338988
$ws.Range("Q5").Value = 0.195229863543
$ws.Range("R5").Value = 1.757068771887
$ws.Range("S5").Value = 0.0101144373489803
$ws.Range("T5").Value = 0.0101144373489803

$ws.Range("I6").Value = 0.02217319039829088
$ws.Range("J6").Value = 0.02217319039829088
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.677453
$ws.Range("N6").Value = 5.032359
$ws.Range("O6").Value = 0.2543757648546719
$ws.Range("P6").Value = 0.2543757648546719
$ws.Range("Q6").Value = 0.108870054606
$ws.Range("R6").Value = 0.9798304914539999
$ws.Range("S6").Value = 0.005640322266833512
$ws.Range("T6").Value = 0.005640322266833512

$ws.Range("I7").Value = 0.02217319039829088
$ws.Range("J7").Value = 0.02217319039829088
$ws.Range("O7").Value = 0.5611569328879042
$ws.Range("P7").Value = 0.5611569328879042
$ws.Range("S7").Value = 0.01244263951624444
$ws.Range("T7").Value = 0.01244263951624444

$ws.Range("I8").Value = 0.02217319039829088
$ws.Range("J8").Value = 0.02217319039829088
$ws.Range("M8").Value = 1.148238333333333
$ws.Range("N8").Value = 3.444715
$ws.Range("O8").Value = 0.174123510034034
$ws.Range("P8").Value = 0.174123510034034
$ws.Range("Q8").Value = 0.07452296430999999
$ws.Range("R8").Value = 0.6707066787899999
$ws.Range("S8").Value = 0.003860873740803349
$ws.Range("T8").Value = 0.003860873740803348

$ws.Range("I9").Value = 0.02217319039829088
$ws.Range("J9").Value = 0.02217319039829088
$ws.Range("M9").Value = 0.06821100000000001
$ws.Range("N9").Value = 0.204633
$ws.Range("O9").Value = 0.01034379222338988
$ws.Range("P9").Value = 0.01034379222338988
$ws.Range("Q9").Value = 0.004427030322000001
$ws.Range("R9").Value = 0.039843272898
$ws.Range("S9").Value = 0.0002293548744095845
$ws.Range("T9").Value = 0.0002293548744095844
